$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells are treated as text so numeric-looking values
# (e.g. "1.00", "0.581") are not silently coerced into numbers.

$changes = @(
    @{Row=2; D='57.882.18'; E='  +3.33%  '},
    @{Row=3; D='3.283.15'; E='  +1.93%  '},
    @{Row=4; E='  +0.00%  '},
    @{Row=5; D='399.45'; E='  +0.85%  '},
    @{Row=6; D='109.68'; E='  -1.11%  '},
    @{Row=7; D='0.581'; E='  +5.62%  '},
    @{Row=8; D='1.00'; E='  +0.02%  '},
    @{Row=9; D='0.625'; E='  +1.15%  '},
    @{Row=10; D='39.50'; E='  +0.85%  '},
    @{Row=11; D='0.0966'; E='  +5.94%  '},
    @{Row=12; E='  +1.40%  '},
    @{Row=13; D='3.791.72'; E='  +1.80%  '},
    @{Row=14; D='8.34'; E='  +3.39%  '},
    @{Row=15; D='19.06'; E='  +0.16%  '},
    @{Row=16; D='3.290.92'; E='  +3.16%  '},
    @{Row=17; D='1.04'; E='  -0.80%  '},
    @{Row=18; D='11.04'; E='  +1.77%  '},
    @{Row=19; D='57.613.85'; E='  +3.17%  '},
    @{Row=20; D='3.33'; E='  +0.10%  '},
    @{Row=21; D='0.0000108'; E='  +5.77%  '},
    @{Row=22; D='12.97'; E='  +0.26%  '},
    @{Row=23; D='300.03'; E='  +0.63%  '},
    @{Row=24; D='74.52'; E='  -1.20%  '},
    @{Row=25; E='  -0.55%  '},
    @{Row=26; D='28.22'; E='  +0.53%  '},
    @{Row=27; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='7.93'; E='  -2.57%  '},
    @{Row=28; B='LEO'; C='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; D='4.41'; E='  +1.00%  '},
    @{Row=29; D='7.39'; E='  -1.37%  '},
    @{Row=30; D='0.170'; E='  -1.45%  '},
    @{Row=32; D='0.113'; E='  +2.22%  '},
    @{Row=33; D='11.29'; E='  +1.66%  '},
    @{Row=34; D='40.78'; E='  +12.87%  '},
    @{Row=35; D='0.0501'; E='  +2.35%  '},
    @{Row=36; E='  +0.95%  '},
    @{Row=37; D='51.72'; E='  +0.63%  '},
    @{Row=38; D='3.15'; E='  +0.28%  '},
    @{Row=39; D='0.998'},
    @{Row=40; D='3.50'; E='  -0.81%  '},
    @{Row=41; D='138.78'; E='  +2.40%  '},
    @{Row=42; D='0.122'; E='  +2.15%  '},
    @{Row=43; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='1.88'; E='  -1.55%  '},
    @{Row=44; B='TheGraph'; C='https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'; D='0.285'; E='  +0.47%  '},
    @{Row=45; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='3.92'; E='  -1.71%  '},
    @{Row=46; B='Celestia'; C='https://coinranking.com/coin/YQcD0lBl7+celestia-tia'; D='16.88'; E='  -2.57%  '},
    @{Row=47; D='22.40'; E='  +1.18%  '},
    @{Row=48; D='2.22'; E='  +4.81%  '},
    @{Row=49; D='2.160.06'; E='  +1.83%  '},
    @{Row=50; D='2.45'; E='  -0.76%  '},
    @{Row=51; D='1.91'; E='  -11.72%  '}
)

foreach ($item in $changes) {
    $r = $item.Row
    if ($item.ContainsKey("B")) { $ws.Cells.Item($r, 2).Value = $item.B }
    if ($item.ContainsKey("C")) { $ws.Cells.Item($r, 3).Value = $item.C }
    if ($item.ContainsKey("D")) {
        $dcell = $ws.Cells.Item($r, 4)
        $dcell.NumberFormat = "@"
        $dcell.Value = $item.D
    }
    if ($item.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $item.E }
}

Write-Host "Applied $($changes.Count) row updates"
